$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Total" row (row 13) that sums the labor cost row (B10:D10)
$ws.Range("A13").Value = "Total"
$ws.Range("B13").Formula = "=SUM(B10:D10)"

# Update the view so it matches the saved selection/scroll state
$ws.Range("B14").Select()
